# Add "Italy" test-data sheet, cloned from "Slovakia", with its own
# market/product values and a one-off larger/borderless font for the
# product code cell.

$wb = $excel.ActiveWorkbook

$slovakia = $wb.Worksheets.Item("Slovakia")

# Clone the Slovakia sheet and drop the copy right after it (end of the
# tab strip), exactly like Excel's own "Move or Copy… > Create a copy".
$slovakia.Copy($null, $slovakia)
$italy = $wb.Worksheets.Item($slovakia.Index + 1)
$italy.Name = "Italy"

# Market name + product code for the new market (product code entered
# first so it lands in the shared-string table ahead of the market name).
$italy.Range("B4").Value = "NGC-3145/T2155"
$italy.Range("B2").Value = "Italy Market"

# The product code was typed in manually with a bigger, borderless font.
$codeCell = $italy.Range("B4")
$codeCell.Font.Size = 12
$codeCell.Font.Color = 0
$codeCell.Borders.LineStyle = 0

# The previously-active Slovakia sheet had its whole grid selected when
# focus moved away (e.g. a stray Ctrl+A before switching tabs).
$slovakia.Activate()
$slovakia.Cells.Select()

# Leave the cursor on the cell that was just edited on the new sheet.
$italy.Activate()
$italy.Range("B2").Select()
